$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds price labels stored as literal text (e.g. "$103"), one row per
# SKU/carrier combination. This block nudges 21 distinct price points by a dollar
# or a few dollars, affecting every row (carrier variant) that shares that price.
# Re-applying a Text number format before writing the value keeps Excel from
# auto-converting the "$NNN" string into a currency-formatted number.

$priceUpdates = @(
    ,@(3, '$102')
    ,@(7, '$102')
    ,@(11, '$102')
    ,@(15, '$102')
    ,@(19, '$102')
    ,@(23, '$140')
    ,@(27, '$140')
    ,@(31, '$140')
    ,@(35, '$140')
    ,@(39, '$140')
    ,@(63, '$179')
    ,@(64, '$192')
    ,@(67, '$179')
    ,@(68, '$192')
    ,@(71, '$179')
    ,@(72, '$192')
    ,@(75, '$179')
    ,@(76, '$192')
    ,@(79, '$179')
    ,@(80, '$192')
    ,@(82, '$187')
    ,@(85, '$187')
    ,@(88, '$187')
    ,@(91, '$187')
    ,@(94, '$187')
    ,@(98, '$281')
    ,@(101, '$281')
    ,@(104, '$281')
    ,@(107, '$281')
    ,@(110, '$281')
    ,@(113, '$312')
    ,@(116, '$312')
    ,@(119, '$312')
    ,@(122, '$312')
    ,@(125, '$312')
    ,@(126, '$336')
    ,@(128, '$389')
    ,@(129, '$336')
    ,@(131, '$389')
    ,@(132, '$336')
    ,@(134, '$389')
    ,@(135, '$336')
    ,@(137, '$389')
    ,@(138, '$336')
    ,@(140, '$389')
    ,@(141, '$408')
    ,@(143, '$408')
    ,@(145, '$408')
    ,@(147, '$408')
    ,@(149, '$408')
    ,@(151, '$419')
    ,@(154, '$419')
    ,@(157, '$419')
    ,@(160, '$419')
    ,@(163, '$419')
    ,@(167, '$540')
    ,@(170, '$540')
    ,@(173, '$540')
    ,@(176, '$540')
    ,@(179, '$540')
    ,@(181, '$542')
    ,@(184, '$542')
    ,@(187, '$542')
    ,@(190, '$542')
    ,@(193, '$542')
    ,@(196, '$607')
    ,@(199, '$607')
    ,@(202, '$607')
    ,@(205, '$607')
    ,@(208, '$607')
    ,@(211, '$773')
    ,@(212, '$884')
    ,@(213, '$962')
    ,@(214, '$773')
    ,@(215, '$884')
    ,@(216, '$962')
    ,@(217, '$773')
    ,@(218, '$884')
    ,@(219, '$962')
    ,@(220, '$773')
    ,@(221, '$884')
    ,@(222, '$962')
    ,@(223, '$773')
    ,@(224, '$884')
    ,@(225, '$962')
    ,@(226, '$860')
    ,@(227, '$948')
    ,@(228, '$1050')
    ,@(229, '$860')
    ,@(230, '$948')
    ,@(231, '$1050')
    ,@(232, '$860')
    ,@(233, '$948')
    ,@(234, '$1050')
    ,@(235, '$860')
    ,@(236, '$948')
    ,@(237, '$1050')
    ,@(238, '$860')
    ,@(239, '$948')
    ,@(240, '$1050')
    ,@(241, '$343')
    ,@(244, '$343')
    ,@(247, '$343')
    ,@(250, '$343')
    ,@(253, '$343')
)

foreach ($update in $priceUpdates) {
    $rowNum = $update[0]
    $newValue = $update[1]
    $cell = $ws.Cells.Item($rowNum, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
}
